$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell A16 content from "asd" to "Test"
$ws.Range("A16").Value = "Test"
